$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure these cells keep their text formatting (values like "2,689" and counts
# are stored as text in this sheet, not numbers) before writing new values.
$ws.Range("B1447:E1478").NumberFormat = "@"

$ws.Range("D1447").Value = "3,471"
$ws.Range("E1447").Value = "3"
$ws.Range("C1448").Value = "Minnesota Vikings"
$ws.Range("D1448").Value = "2,924"
$ws.Range("E1448").Value = "2"
$ws.Range("C1449").Value = "Houston Texans"
$ws.Range("D1449").Value = "3,020"
$ws.Range("E1449").Value = "1"
$ws.Range("C1450").Value = "Tampa Bay Buccaneers"
$ws.Range("D1450").Value = "3,073"
$ws.Range("E1450").Value = "3"
$ws.Range("C1451").Value = "Indianapolis Colts"
$ws.Range("D1451").Value = "2,737"
$ws.Range("E1451").Value = "3"
$ws.Range("C1452").Value = "Seattle Seahawks"
$ws.Range("D1452").Value = "2,853"
$ws.Range("E1452").Value = "3"
$ws.Range("D1453").Value = "2,809"
$ws.Range("D1454").Value = "2,452"
$ws.Range("E1454").Value = "4"
$ws.Range("C1455").Value = "New York Jets"
$ws.Range("D1455").Value = "2,704"
$ws.Range("E1455").Value = "2"
$ws.Range("C1456").Value = "Atlanta Falcons"
$ws.Range("D1456").Value = "2,758"
$ws.Range("C1457").Value = "Tennessee Titans"
$ws.Range("D1457").Value = "2,736"
$ws.Range("C1458").Value = "Detroit Lions"
$ws.Range("D1458").Value = "3,024"
$ws.Range("E1458").Value = "4"
$ws.Range("C1459").Value = "Los Angeles Rams"
$ws.Range("D1459").Value = "2,527"
$ws.Range("E1459").Value = "4"
$ws.Range("D1460").Value = "2,554"
$ws.Range("C1461").Value = "Dallas Cowboys"
$ws.Range("D1461").Value = "2,372"
$ws.Range("B1462").Value = "16"
$ws.Range("C1462").Value = "New Orleans Saints"
$ws.Range("D1462").Value = "2,513"
$ws.Range("B1463").Value = "16"
$ws.Range("C1463").Value = "Washington Redskins"
$ws.Range("D1463").Value = "2,492"
$ws.Range("B1464").Value = "18"
$ws.Range("C1464").Value = "Los Angeles Chargers"
$ws.Range("D1464").Value = "2,417"
$ws.Range("E1464").Value = "1"
$ws.Range("B1465").Value = "18"
$ws.Range("C1465").Value = "Pittsburgh Steelers"
$ws.Range("D1465").Value = "2,455"
$ws.Range("E1465").Value = "2"
$ws.Range("C1466").Value = "Oakland Raiders"
$ws.Range("D1466").Value = "2,817"
$ws.Range("E1466").Value = "1"
$ws.Range("B1467").Value = "20"
$ws.Range("C1467").Value = "Philadelphia Eagles"
$ws.Range("D1467").Value = "2,540"
$ws.Range("E1467").Value = "1"
$ws.Range("B1468").Value = "22"
$ws.Range("C1468").Value = "Buffalo Bills"
$ws.Range("D1468").Value = "2,167"
$ws.Range("E1468").Value = "7"
$ws.Range("B1469").Value = "22"
$ws.Range("C1469").Value = "Green Bay Packers"
$ws.Range("D1469").Value = "2,756"
$ws.Range("E1469").Value = "2"
$ws.Range("C1470").Value = "Denver Broncos"
$ws.Range("D1470").Value = "2,232"
$ws.Range("E1470").Value = "3"
$ws.Range("C1471").Value = "New York Giants"
$ws.Range("D1471").Value = "2,755"
$ws.Range("E1471").Value = "2"
$ws.Range("C1472").Value = "Baltimore Ravens"
$ws.Range("D1472").Value = "2,555"
$ws.Range("E1472").Value = "6"
$ws.Range("C1473").Value = "Jacksonville Jaguars"
$ws.Range("D1473").Value = "2,408"
$ws.Range("E1473").Value = "2"
$ws.Range("C1474").Value = "Cleveland Browns"
$ws.Range("D1474").Value = "2,379"
$ws.Range("E1474").Value = "2"
$ws.Range("B1475").Value = "29"
$ws.Range("C1475").Value = "Miami Dolphins"
$ws.Range("D1475").Value = "2,545"
$ws.Range("E1475").Value = "0"
$ws.Range("B1476").Value = "30"
$ws.Range("C1476").Value = "Cincinnati Bengals"
$ws.Range("D1476").Value = "2,661"
$ws.Range("E1476").Value = "3"
$ws.Range("C1477").Value = "New England Patriots"
$ws.Range("D1477").Value = "1,780"
$ws.Range("E1477").Value = "2"
$ws.Range("D1478").Value = "1,735"
$ws.Range("E1478").Value = "4"
